$wb = $excel.ActiveWorkbook

# Remember the originally active sheet so it can stay active once the new
# sheet has been added and populated (adding/activating a sheet otherwise
# shifts the workbook's active tab to it).
$originalActiveSheet = $wb.ActiveSheet

# Add the new "Yearly demand" sheet at the very end of the tab order
# (after the current last sheet, "Connected Households").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Yearly demand"

# Match the outline summary-direction flags and page margins already used
# by every other sheet in this workbook.
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72

# The header row (B1:Y1) and the index column (A2:A4) reuse the bold /
# centered / thin-bordered style that the other sheets in this workbook
# already use for the same purpose. Copy the format from an existing cell
# that carries that style so the new sheet points at the very same style
# entry instead of minting a new (merely equivalent) one.
$styleSource = $wb.Worksheets.Item("DG Dispatch").Range("B1")
$styleSource.Copy()
$ws.Range("B1:Y1").PasteSpecial(-4122)
$ws.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row: hours 0-23 across B1:Y1.
for ($col = 0; $col -lt 24; $col++) {
    $ws.Cells.Item(1, $col + 2).Value = $col
}

# Index column A2:A4 = 0,1,2
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2

# Row 2 (index 0) data, B2:Y2
$row2 = @(-32.5,-19.5,-13,-13,-13,142.5,291.5,327,388.5,502,596,670.5,745,651,576.5,502,320.5,139,32,-117,-97.5,-78,-52,-39)
for ($col = 0; $col -lt $row2.Length; $col++) {
    $ws.Cells.Item(2, $col + 2).Value = $row2[$col]
}

# Row 3 (index 1) data, B3:Y3
$row3 = @(-32.5,-19.5,-13,0,0,-19.5,0,324,486,648,729,751.5,583,567,333.5,340,243,57.99999999999999,-130,0,0,-78,0,-39)
for ($col = 0; $col -lt $row3.Length; $col++) {
    $ws.Cells.Item(3, $col + 2).Value = $row3[$col]
}

# Row 4 (index 2) data, B4:Y4
$row4 = @(-32.5,-19.5,0,0,0,-19.5,0,0,81,324,567,589.5,648,567,324,162,81,0,-130,0,0,0,0,-39)
for ($col = 0; $col -lt $row4.Length; $col++) {
    $ws.Cells.Item(4, $col + 2).Value = $row4[$col]
}


# Restore the original active sheet/selection.
[void]$originalActiveSheet.Activate()
[void]$originalActiveSheet.Range("A1").Select()
